$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The slide already has three rectangles (ids 4, 5, 6). The real PowerPoint
# id allocator would hand a freshly-created shape id 7 (next unused id after
# the existing 1,4,5,6), but this runtime's allocator starts handing out the
# low ids 2 and 3 first. Burn through those two with throwaway shapes (which
# we immediately delete) so the shape we actually care about lands on id 7,
# matching "Rectangle 6" / id="7" from the target deck.
$tmp1 = $s.Shapes.AddShape(1, 0, 0, 10, 10)
$tmp2 = $s.Shapes.AddShape(1, 0, 0, 10, 10)
$tmp1.Delete()
$tmp2.Delete()

# Duplicate the "CSC-315" rectangle (id 5) - it already has the right
# accent2 fill / style / size - to create the new "CSC-450" rectangle below
# the empty "Rectangle 5" (id 6), then move it into place.
$src = $s.Shapes.Item(2)
$new = $src.Duplicate()
$new.Name = "Rectangle 6"

# Shape.Left/.Top are in points; the target position is EMU (5836666, 3993810).
# Add a tiny epsilon so the points->EMU round trip inside the host lands on
# the exact target EMU instead of being floored one unit short.
$new.Left = 5836666 / 12700 + 0.00002
$new.Top = 3993810 / 12700 + 0.00002

# Target text is two runs: "CSC" then "-450".
$tr = $new.TextFrame.TextRange
$tr.Text = "CSC"
$tr.InsertAfter("-450") | Out-Null
